$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 contains three "blog" widgets (C7, E7, I7) whose "ser:" numbers
# are bumped by one because a new blog post was uploaded.

# C7: blog ser 99 -> ser 100
$ws.Range("C7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 100"

# E7: blog ser 98 -> ser 99
$ws.Range("E7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 99"

# I7: blog ser 97 -> ser 98
$ws.Range("I7").Value = "type: blog`nwidth: 2`nheight: 1`nser: 98"
